$wb = $excel.ActiveWorkbook

# --- Fix the "Enrolment Statistics" -> "Enrollment Statistics" spelling everywhere it is reachable ---
$oldName = "Enrolment Statistics"
$newName = "Enrollment Statistics"

# Locate the worksheet by its current (misspelled) name and rename it.
$wsStats = $wb.Worksheets.Item($oldName)
$wsStats.Name = $newName

# The pie chart embedded on that worksheet still has series formulas pointing at the
# old sheet name (renaming the sheet does not automatically rewrite chart series
# formulas), so update them explicitly.
for ($c = 1; $c -le $wsStats.ChartObjects().Count; $c++) {
    $co = $wsStats.ChartObjects().Item($c)
    $chart = $co.Chart
    $series = $chart.SeriesCollection()
    for ($i = 1; $i -le $series.Count; $i++) {
        $s = $series.Item($i)
        $formula = $s.Formula
        if ($formula -like "*$oldName*") {
            $s.Formula = $formula.Replace($oldName, $newName)
        }
    }
}

# Make the renamed "Enrollment Statistics" sheet the active / selected tab (it becomes
# the last-used sheet after the rename and related tweaks).
$wsStats.Activate()
